$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("900-1145PM")
$ws.Name = "900-1145 PM"
